$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data

$ws.Range("D2").Value = "46.643.16"
$ws.Range("E2").Value = "  +6.22%  "

$ws.Range("D3").Value = "2.298.50"
$ws.Range("E3").Value = "  +3.57%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.52"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.86"
$ws.Range("E6").Value = "  +11.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  +2.04%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +6.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.66"
$ws.Range("E10").Value = "  +11.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  +2.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.37"
$ws.Range("E12").Value = "  +6.27%  "

$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").Value = "2.646.78"
$ws.Range("E14").Value = "  +3.53%  "

$ws.Range("D15").Value = "2.295.24"
$ws.Range("E15").Value = "  +3.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.82"
$ws.Range("E16").Value = "  +3.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.812"
$ws.Range("E17").Value = "  +4.85%  "

$ws.Range("D18").Value = "46.600.18"
$ws.Range("E18").Value = "  +6.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").Value = "  +11.41%  "

$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  +3.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.99"
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.44"
$ws.Range("E22").Value = "  +3.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.91"
$ws.Range("E23").Value = "  +5.50%  "

$ws.Range("E24").Value = "  +3.12%  "

$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").Value = "  +3.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.85"
$ws.Range("E27").Value = "  +11.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.78"
$ws.Range("E29").Value = "  +4.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.01"
$ws.Range("E30").Value = "  +4.56%  "

$ws.Range("E31").Value = "  +12.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.65"
$ws.Range("E32").Value = "  +4.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.93"
$ws.Range("E33").Value = "  -4.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0796"
$ws.Range("E34").Value = "  +5.01%  "

$ws.Range("E35").Value = "  +16.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  +12.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  +6.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.98"
$ws.Range("E39").Value = "  +20.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.99"
$ws.Range("E40").Value = "  +10.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.36"
$ws.Range("E41").Value = "  +6.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0302"
$ws.Range("E42").Value = "  +1.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +9.18%  "

$ws.Range("D45").Value = "1.810.50"
$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.92"
$ws.Range("E46").Value = "  +21.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.195"
$ws.Range("E47").Value = "  +6.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.27"
$ws.Range("E48").Value = "  +9.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.91"
$ws.Range("E49").Value = "  +6.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.95"
$ws.Range("E50").Value = "  +1.73%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.71"
$ws.Range("E51").Value = "  +5.28%  "
